$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray spaces around the "+" in these two group labels
$ws.Range("B4").Value = "3013+3015"
$ws.Range("B6").Value = "3016+3019"

# Leave the selection on B6, matching where the edit was made
$ws.Range("B6").Select()
